$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" (D) column holds plain numeric-looking strings (e.g. "228.88")
# that must stay as text, matching the original inline-string cells, instead
# of being auto-converted to a Number by Excel (which would also drop
# formatting like trailing zeros, e.g. "1.00" -> 1). Force each touched D
# cell to Text format before writing its new value.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '38.072.38'
$ws.Range('E2').Value = '  +0.35%  '
$ws.Range('D3').Value = '2.040.34'
$ws.Range('E3').Value = '  -0.57%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '228.88'
$ws.Range('E5').Value = '  -0.41%  '
$ws.Range('E6').Value = '  -0.51%  '
$ws.Range('D7').Value = '60.82'
$ws.Range('E7').Value = '  +3.83%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('D10').Value = '0.0821'
$ws.Range('E10').Value = '  +1.36%  '
$ws.Range('D11').Value = '0.104'
$ws.Range('E11').Value = '  +0.49%  '
$ws.Range('D12').Value = '14.69'
$ws.Range('E12').Value = '  +0.30%  '
$ws.Range('D13').Value = '2.341.62'
$ws.Range('E13').Value = '  -0.58%  '
$ws.Range('D14').Value = '21.43'
$ws.Range('E14').Value = '  +2.84%  '
$ws.Range('D15').Value = '0.769'
$ws.Range('E15').Value = '  +1.88%  '
$ws.Range('D16').Value = '5.19'
$ws.Range('E16').Value = '  -1.59%  '
$ws.Range('D17').Value = '2.041.41'
$ws.Range('E17').Value = '  -0.70%  '
$ws.Range('D18').Value = '37.867.70'
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = '5.98'
$ws.Range('E19').Value = '  -5.03%  '
$ws.Range('B20').Value = 'Litecoin'
$ws.Range('C20').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D20').Value = '69.94'
$ws.Range('E20').Value = '  +0.32%  '
$ws.Range('D21').Value = '0.0₃0827'
$ws.Range('E21').Value = '  -1.18%  '
$ws.Range('D22').Value = '224.99'
$ws.Range('E22').Value = '  +0.28%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('E24').Value = '  +0.28%  '
$ws.Range('E25').Value = '  -0.57%  '
$ws.Range('D26').Value = '9.36'
$ws.Range('E26').Value = '  +0.69%  '
$ws.Range('D27').Value = '167.34'
$ws.Range('E27').Value = '  +0.48%  '
$ws.Range('E28').Value = '  -2.32%  '
$ws.Range('D29').Value = '18.95'
$ws.Range('E29').Value = '  -0.36%  '
$ws.Range('E30').Value = '  -3.09%  '
$ws.Range('E31').Value = '  +0.75%  '
$ws.Range('E32').Value = '  +5.85%  '
$ws.Range('E33').Value = '  -2.15%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.0609'
$ws.Range('E34').Value = '  -0.12%  '
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').Value = '4.54'
$ws.Range('E35').Value = '  -1.15%  '
$ws.Range('E36').Value = '  +6.58%  '
$ws.Range('E37').Value = '  -1.81%  '
$ws.Range('E38').Value = '  +0.32%  '
$ws.Range('E39').Value = '  -0.11%  '
$ws.Range('D40').Value = '17.92'
$ws.Range('E40').Value = '  +8.30%  '
$ws.Range('D41').Value = '1.528.49'
$ws.Range('E41').Value = '  +1.96%  '
$ws.Range('E42').Value = '  +0.75%  '
$ws.Range('D43').Value = '96.46'
$ws.Range('E43').Value = '  -0.82%  '
$ws.Range('E44').Value = '  -2.67%  '
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('E46').Value = '  -2.50%  '
$ws.Range('D47').Value = '4.06'
$ws.Range('E47').Value = '  -1.62%  '
$ws.Range('D48').Value = '1.01'
$ws.Range('E48').Value = '  -0.20%  '
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('D50').Value = '7.13'
$ws.Range('E50').Value = '  +0.46%  '
$ws.Range('D51').Value = '2.231.60'
$ws.Range('E51').Value = '  -0.56%  '
